# ---------------------------------------------------------------------------
# Adds an "Add-Users" worksheet to Users.xlsx and reshapes "Get-Users" so the
# previously nested address/geo/company JSON fields are flattened into plain
# columns (street, suite, city, zipcode, lat, lng, ... companyName, ...).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Flatten "Get-Users": drop the wrapper-label columns ("address",
#        "geo", "company") so the nested fields sit directly under the
#        top-level header row. Deleting whole columns (rather than rewriting
#        cell-by-cell) lets any per-cell formatting (e.g. the quirky font on
#        the "suite" header) ride along with the shift, just like it would
#        in the live UI.
$ws.Columns.Item(15).Delete()   # O: company
$ws.Columns.Item(10).Delete()   # J: geo
$ws.Columns.Item(5).Delete()    # E: address

# --- 2. Fix the "company name" header: it used to be a duplicate of the
#        "name" header (column P); give it its own proper label.
$ws.Range("M1").Value = "companyName"

# --- 3. Populate the newly-exposed street/suite/city/zipcode columns and the
#        lat/lng/company columns with the flattened user records (still the
#        same two JSONPlaceholder users as before, just fully expanded).
$ws.Range("E2").Value = "Kulas Light"
$ws.Range("F2").Value = "Apt. 556"
$ws.Range("G2").Value = "Gwenborough"
$ws.Range("H2").Value = "92998-3874"
$ws.Range("M2").Value = "Romaguera-Crona"
$ws.Range("N2").Value = "Multi-layered client-server neural-net"
$ws.Range("O2").Value = "harness real-time e-markets"

$ws.Range("E3").Value = "Victor Plains"
$ws.Range("F3").Value = "Suite 879"
$ws.Range("G3").Value = "Wisokyburgh"
$ws.Range("H3").Value = "90566-7771"
$ws.Range("M3").Value = "Deckow-Crist"
$ws.Range("N3").Value = "Proactive didactic contingency"
$ws.Range("O3").Value = "synergize scalable supply-chains"

# lat/lng ("I"/"J") are latitude/longitude text that looks numeric, so force
# a text format first (otherwise "-37.3159" silently becomes the float
# -37.315899999999999) - mirrors the "Format Cells > Text" + ignore-error
# flag Excel shows for these two columns in the real workbook.
$ws.Range("I1:J3").NumberFormat = "@"
$ws.Range("I2").Value = "-37.3159"
$ws.Range("J2").Value = "81.1496"
$ws.Range("I3").Value = "-43.9509"
$ws.Range("J3").Value = "-34.4618"

# --- 4. A couple of columns were manually widened once the flattened data
#        was in place.
[void]$ws.Columns.Item(4).ColumnWidth  # D: email (no-op placeholder removed below)
$ws.Columns.Item(4).ColumnWidth = 16.666666666666668   # D: email  -> 17.5
$ws.Columns.Item(9).ColumnWidth = 14.166666666666666   # I: lat    -> 15
$ws.Columns.Item(14).ColumnWidth = 31.666666666666668  # N: catchPhrase -> 32.5

# --- 5. Re-select the cell the author last had active, and make sure
#        Get-Users keeps being the visible/active tab.
$ws.Range("M11").Select()
$ws.Activate()

# ---------------------------------------------------------------------------
# 6. Add the "Add-Users" worksheet right after "Get-Users".
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Add-Users"

$ws2.Range("A1").Value = "name"
$ws2.Range("B1").Value = "username"
$ws2.Range("C1").Value = "email"
$ws2.Range("D1").Value = "street"
$ws2.Range("E1").Value = "suite"
$ws2.Range("F1").Value = "city"
$ws2.Range("G1").Value = "zipcode"
$ws2.Range("H1").Value = "lat"
$ws2.Range("I1").Value = "lng"
$ws2.Range("J1").Value = "phone"
$ws2.Range("K1").Value = "website"
$ws2.Range("L1").Value = "companyName"
$ws2.Range("M1").Value = "catchPhrase"
$ws2.Range("N1").Value = "bs"

$ws2.Range("A2").Value = "Cristiano Ronaldo"
$ws2.Range("B2").Value = "Cr7"
$ws2.Range("C2").Value = "cr7@april.biz"
$ws2.Range("D2").Value = "Buenos Street 01"
$ws2.Range("E2").Value = "Apt. 001"
$ws2.Range("F2").Value = "Paolo 001"
$ws2.Range("G2").Value = 1232575
$ws2.Range("H2").Value = -390192
$ws2.Range("I2").Value = 109780
$ws2.Range("J2").Value = "1-770-736-8031 x56442"
$ws2.Range("K2").Value = "cr7.org"
$ws2.Range("L2").Value = "cr7 corp"
$ws2.Range("M2").Value = "best player"
$ws2.Range("N2").Value = "best player 001"

$ws2.Range("A3").Value = "Lionel Messi"
$ws2.Range("B3").Value = "LM"
$ws2.Range("C3").Value = "messi@melissa.tv"
$ws2.Range("D3").Value = "Buenos Street 02"
$ws2.Range("E3").Value = "Apt. 002"
$ws2.Range("F3").Value = "Paolo 002"
$ws2.Range("G3").Value = 7047346
$ws2.Range("H3").Value = -390194
$ws2.Range("I3").Value = 109710
$ws2.Range("J3").Value = "010-692-6593 x09125"
$ws2.Range("K3").Value = "lm.org"
$ws2.Range("L3").Value = "lm corp"
$ws2.Range("M3").Value = "best match"
$ws2.Range("N3").Value = "best match 001"

# Hyperlinks on the e-mail column (carried over from the Get-Users cells they
# were copied from, so the link target is still the original addresses even
# though the visible text is the new placeholder e-mail).
$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:Sincere@april.biz", "", "", "Sincere@april.biz")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "mailto:Shanna@melissa.tv", "", "", "Shanna@melissa.tv")

# Column widths that were carried over / re-fitted on the new sheet.
$ws2.Columns.Item(1).ColumnWidth = 14.666666666666666   # A: name        -> 15.5
$ws2.Columns.Item(3).ColumnWidth = 14.666666666666666   # C: email       -> 15.5
$ws2.Columns.Item(4).ColumnWidth = 14.333333333333332   # D: street      -> 15.1640625
$ws2.Columns.Item(10).ColumnWidth = 19.833333333333332  # J: phone       -> 20.6640625
$ws2.Columns.Item(12).ColumnWidth = 12.666666666666666  # L: companyName -> 13.5
$ws2.Columns.Item(14).ColumnWidth = 12.999999999999998  # N: bs          -> 13.83203125

$ws2.Range("M14").Select()
